# This edit corresponds to a fixture regeneration (POI packaging fix / POI
# version upgrade) that only changed *how* the existing OOXML was
# serialized (XML attribute ordering inside word/document.xml and
# word/styles.xml, e.g. <w:pgSz w:w=".." w:h=".."/> -> <w:pgSz w:h=".." w:w=".."/>).
# Every changed line in the source diff carries the exact same element
# name, attribute names and attribute values as before -- only the
# attribute order differs. There is no visible-content, formatting,
# structural, or style-value change to apply via the Word object model:
# the document's text, styles, numbering, sections and properties are
# already identical to the target state.
$d = $word.ActiveDocument
